# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) for most rows, and for rows 43-45 also
# updates Coin (B) and Link (C) because three coins shifted rank order
# (Stacks/Hedera/FirstDigitalUSD -> FirstDigitalUSD/Stacks/Hedera).
#
# Several Price values are plain decimals (e.g. "7.40", "1.60", "0.0890")
# that Excel would otherwise auto-coerce to numbers and strip the
# significant trailing zero (7.40 -> 7.4). Forcing the cell to Text via
# NumberFormat "@" before the write, then clearing the format afterwards,
# keeps the literal string without leaving a stray number-format style
# behind (matches the original file, where these are plain, unstyled text
# cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.099.38"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "2.829.59"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.61%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -5.22%  "
$ws.Range("D9").Value = "2.829.41"
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.92"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("D13").Value = "3.332.64"
$ws.Range("E13").Value = "  -3.77%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "59.235.79"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.95%  "
$ws.Range("D17").Value = "2.847.83"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("E18").Value = "  -6.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.69%  "
$ws.Range("E22").Value = "  -5.51%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("E26").Value = "  -6.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.172"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -7.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.67%  "
$ws.Range("D30").Value = "0.0₃0801"
$ws.Range("E30").Value = "  -10.35%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.60"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.05"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.16"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.35"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("E37").Value = "  -11.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").Value = "2.213.01"
$ws.Range("E40").Value = "  -7.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.628"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.47%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -9.92%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0558"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -9.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0890"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -10.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.56%  "
